$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6668.4287
$ws.Range("I43").Value = 6668.4287
$ws.Range("K43").Value = 6668.4287
$ws.Range("M43").Value = -6599.4287
$ws.Range("H51").Value = 5571.684
$ws.Range("I51").Value = 4569.4287
$ws.Range("J51").Value = 6156.3335
$ws.Range("K51").Value = 4569.4287
$ws.Range("L51").Value = 6156.3335
$ws.Range("M51").Value = -4085.4287
$ws.Range("N51").Value = -7124.3335
$ws.Range("H98").Value = 1841.2258
$ws.Range("I98").Value = 1645.5555
$ws.Range("K98").Value = 1645.5555
$ws.Range("M98").Value = -147.5554999999999
$ws.Range("H116").Value = 3569.6
$ws.Range("I116").Value = 3775.05
$ws.Range("J116").Value = 3158.7
$ws.Range("K116").Value = 3775.05
$ws.Range("L116").Value = 3158.7
$ws.Range("M116").Value = -333.0500000000002
$ws.Range("N116").Value = -10042.7
$ws.Range("H122").Value = 1841.2258
$ws.Range("I122").Value = 1645.5555
$ws.Range("K122").Value = 4936.666499999999
$ws.Range("M122").Value = -2486.666499999999
$ws.Range("H131").Value = 5369.9287
$ws.Range("I131").Value = 1475
$ws.Range("K131").Value = 4425
$ws.Range("M131").Value = 615
$ws.Range("H137").Value = 2455.1035
$ws.Range("I137").Value = 2115.5
$ws.Range("J137").Value = 3209.7778
$ws.Range("K137").Value = 6346.5
$ws.Range("L137").Value = 9629.3334
$ws.Range("M137").Value = -3796.5
$ws.Range("N137").Value = -14729.3334
$ws.Range("H138").Value = 2369.8438
$ws.Range("I138").Value = 1235.4242
$ws.Range("J138").Value = 3577.4517
$ws.Range("K138").Value = 3706.2726
$ws.Range("L138").Value = 10732.3551
$ws.Range("M138").Value = 1433.7274
$ws.Range("N138").Value = -21012.3551

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2103671.8
$ws.Range("I2").Value = 2452783.8
$ws.Range("J2").Value = 8999
$ws.Range("K2").Value = 2452783.8
$ws.Range("L2").Value = 8999
$ws.Range("M2").Value = -2452670.8
$ws.Range("N2").Value = -9225
$ws.Range("H32").Value = 2437.8293
$ws.Range("I32").Value = 2256.7742
$ws.Range("J32").Value = 2999.1
$ws.Range("K32").Value = 2256.7742
$ws.Range("L32").Value = 2999.1
$ws.Range("M32").Value = -1969.7742
$ws.Range("N32").Value = -3573.1
$ws.Range("H110").Value = 73771.92999999999
$ws.Range("I110").Value = 85233.914
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 85233.914
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -83188.914
$ws.Range("N110").Value = -9090
$ws.Range("H116").Value = 2103671.8
$ws.Range("I116").Value = 2452783.8
$ws.Range("J116").Value = 8999
$ws.Range("K116").Value = 2452783.8
$ws.Range("L116").Value = 8999
$ws.Range("M116").Value = -2450489.8
$ws.Range("N116").Value = -13587
$ws.Range("H132").Value = 2944516.2
$ws.Range("I132").Value = 3451576.5
$ws.Range("J132").Value = 3566.8
$ws.Range("K132").Value = 10354729.5
$ws.Range("L132").Value = 10700.4
$ws.Range("M132").Value = -10352199.5
$ws.Range("N132").Value = -15760.4

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2103671.8
$ws.Range("I3").Value = 2452783.8
$ws.Range("J3").Value = 8999
$ws.Range("K3").Value = 2452783.8
$ws.Range("L3").Value = 8999
$ws.Range("M3").Value = -2452669.8
$ws.Range("N3").Value = -9227
$ws.Range("H94").Value = 15350.625
$ws.Range("I94").Value = 15558
$ws.Range("J94").Value = 13899
$ws.Range("K94").Value = 15558
$ws.Range("L94").Value = 13899
$ws.Range("M94").Value = -15107
$ws.Range("N94").Value = -14801
$ws.Range("H109").Value = 20000
$ws.Range("I109").Value = 20000
$ws.Range("K109").Value = 20000
$ws.Range("M109").Value = -18613

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10351.519
$ws.Range("I31").Value = 7721.846
$ws.Range("J31").Value = 12793.357
$ws.Range("K31").Value = 7721.846
$ws.Range("L31").Value = 12793.357
$ws.Range("M31").Value = -7426.846
$ws.Range("N31").Value = -13383.357
$ws.Range("H34").Value = 10351.519
$ws.Range("I34").Value = 7721.846
$ws.Range("J34").Value = 12793.357
$ws.Range("K34").Value = 7721.846
$ws.Range("L34").Value = 12793.357
$ws.Range("M34").Value = -7519.846
$ws.Range("N34").Value = -13197.357
$ws.Range("H99").Value = 3585.75
$ws.Range("I99").Value = 3585.75
$ws.Range("K99").Value = 3585.75
$ws.Range("M99").Value = -2087.75
$ws.Range("H122").Value = 2551.1292
$ws.Range("I122").Value = 2551.1292
$ws.Range("K122").Value = 7653.3876
$ws.Range("M122").Value = -5203.3876
$ws.Range("H126").Value = 3585.75
$ws.Range("I126").Value = 3585.75
$ws.Range("K126").Value = 10757.25
$ws.Range("M126").Value = -8287.25
$ws.Range("H132").Value = 17858774
$ws.Range("I132").Value = 18869584
$ws.Range("K132").Value = 56608752
$ws.Range("M132").Value = -56606222

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 436.6
$ws.Range("I14").Value = 436.6
$ws.Range("K14").Value = 1309.8
$ws.Range("M14").Value = -1136.8

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11473.667
$ws.Range("I102").Value = 1894.7142
$ws.Range("J102").Value = 45000
$ws.Range("K102").Value = 1894.7142
$ws.Range("L102").Value = 45000
$ws.Range("M102").Value = -272.7141999999999
$ws.Range("N102").Value = -48244

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3976.3333
$ws.Range("I7").Value = 3938.5386
$ws.Range("K7").Value = 3938.5386
$ws.Range("M7").Value = -3826.5386
$ws.Range("H46").Value = 1053.5714
$ws.Range("I46").Value = 1143.3334
$ws.Range("J46").Value = 986.25
$ws.Range("K46").Value = 1143.3334
$ws.Range("L46").Value = 986.25
$ws.Range("M46").Value = -955.3334
$ws.Range("N46").Value = -1362.25
$ws.Range("H100").Value = 19962370
$ws.Range("I100").Value = 28515528
$ws.Range("K100").Value = 28515528
$ws.Range("M100").Value = -28514987
$ws.Range("H126").Value = 3976.3333
$ws.Range("I126").Value = 3938.5386
$ws.Range("K126").Value = 11815.6158
$ws.Range("M126").Value = -9345.6158

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15676.286
$ws.Range("I74").Value = 15664
$ws.Range("J74").Value = 15685.5
$ws.Range("K74").Value = 15664
$ws.Range("L74").Value = 15685.5
$ws.Range("M74").Value = -14728
$ws.Range("N74").Value = -17557.5
$ws.Range("H77").Value = 15676.286
$ws.Range("I77").Value = 15664
$ws.Range("J77").Value = 15685.5
$ws.Range("K77").Value = 46992
$ws.Range("L77").Value = 47056.5
$ws.Range("M77").Value = -42312
$ws.Range("N77").Value = -56416.5
$ws.Range("H124").Value = 49999
$ws.Range("J124").Value = 49999
$ws.Range("L124").Value = 49999
$ws.Range("N124").Value = -59819
